# Initial commit of powerplant IO types
#
# Updates the Ohm symbol row to also carry the ASCII alias "ohm", and
# appends three new "miscellaneous" quantity rows describing the
# power-related units VA (volt-ampere), kVA (kilovolt-ampere) and
# var (volt-ampere reactive).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Ohm symbol used to be the bare letter "O" - extend it so it also
# matches the spelled-out alias "ohm".
$ws.Range("B19").Value = "O,ohm"

# New row: volt-ampere (apparent power)
$ws.Range("A159").Value = "miscellaneous"
$ws.Range("B159").Value = "VA"
$ws.Range("C159").Value = "volt ampere"
$ws.Range("D159").Value = -3
$ws.Range("E159").Value = 2
$ws.Range("F159").Value = 1
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 0
$ws.Range("I159").Value = 0
$ws.Range("J159").Value = 0
$ws.Range("K159").Value = 0
$ws.Range("L159").Value = 0
$ws.Range("M159").Value = 1
$ws.Range("N159").Value = "1 volt ampere (measures apparent power)"

# New row: kilovolt-ampere
$ws.Range("A160").Value = "miscellaneous"
$ws.Range("B160").Value = "kVA"
$ws.Range("C160").Value = "kilovolt ampere"
$ws.Range("D160").Value = -3
$ws.Range("E160").Value = 2
$ws.Range("F160").Value = 1
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0
$ws.Range("I160").Value = 0
$ws.Range("J160").Value = 0
$ws.Range("K160").Value = 0
$ws.Range("L160").Value = 0
$ws.Range("M160").Value = 1000
$ws.Range("N160").Value = "VA to kVA"

# New row: volt-ampere reactive (reactive power)
$ws.Range("A161").Value = "miscellaneous"
$ws.Range("B161").Value = "var"
$ws.Range("C161").Value = "volt ampere reactive"
$ws.Range("D161").Value = -3
$ws.Range("E161").Value = 2
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0
$ws.Range("I161").Value = 0
$ws.Range("J161").Value = 0
$ws.Range("K161").Value = 0
$ws.Range("L161").Value = 0
$ws.Range("M161").Value = 1
$ws.Range("N161").Value = "1 volt ampere-reactive (measures reactive power)"

# Leave the view the way the author left it when they saved: scrolled back
# to the top of the frozen pane, with the cursor sitting on G19.
$ws.Activate()
$ws.Range("G19").Select()
